$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 571
$ws.Range("I2").Value = 1536
$ws.Range("J2").Value = 6446
$ws.Range("K2").Value = 35
$ws.Range("L2").Value = 1802
$ws.Range("M2").Value = 109
$ws.Range("N2").Value = 1131
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 24
$ws.Range("Q2").Value = 8
$ws.Range("R2").Value = 63
$ws.Range("S2").Value = 699
$ws.Range("T2").Value = 1175
$ws.Range("U2").Value = 74
$ws.Range("V2").Value = 9842
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 9856
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 162
$ws.Range("AA2").Value = 55
